$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("capacities")
$ws.Activate()

# Replace "water" with "ROR" for the fuel column entries in rows 23-29
$range = $ws.Range("A23:A29")
$range.Value = "ROR"
[void]$range.Select()
